$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextValue "D2" "71.322.27"
$ws.Range("E2").Value = "  +0.46%  "
Set-TextValue "D3" "3.810.43"
$ws.Range("E3").Value = "  -1.03%  "
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue "D5" "700.97"
$ws.Range("E5").Value = "  -0.85%  "
Set-TextValue "D6" "171.11"
$ws.Range("E6").Value = "  -0.97%  "
Set-TextValue "D7" "3.811.58"
$ws.Range("E7").Value = "  -0.94%  "
Set-TextValue "D8" "0.999"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("E10").Value = "  -1.79%  "
Set-TextValue "D11" "7.49"
$ws.Range("E11").Value = "  +1.79%  "
Set-TextValue "D12" "0.485"
$ws.Range("E12").Value = "  +5.74%  "
$ws.Range("E13").Value = "  -2.12%  "
Set-TextValue "D14" "36.02"
$ws.Range("E14").Value = "  -1.92%  "
Set-TextValue "D15" "4.454.39"
$ws.Range("E15").Value = "  -0.99%  "
Set-TextValue "D16" "3.812.46"
$ws.Range("E16").Value = "  -0.57%  "
Set-TextValue "D17" "71.359.85"
$ws.Range("E17").Value = "  +0.48%  "
Set-TextValue "D18" "7.23"
$ws.Range("E18").Value = "  +0.26%  "
Set-TextValue "D19" "17.54"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("E20").Value = "  -0.35%  "
Set-TextValue "D21" "514.52"
$ws.Range("E21").Value = "  +3.53%  "
$ws.Range("E22").Value = "  -1.65%  "
Set-TextValue "D23" "0.715"
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("E24").Value = "  -1.74%  "
$ws.Range("E25").Value = "  -3.21%  "
Set-TextValue "D26" "12.65"
$ws.Range("E26").Value = "  +3.68%  "
Set-TextValue "D27" "3.959.36"
$ws.Range("E27").Value = "  -1.05%  "
$ws.Range("E28").Value = "  -2.64%  "
$ws.Range("E29").Value = "  +0.11%  "
Set-TextValue "D30" "2.02"
$ws.Range("E30").Value = "  -4.26%  "
Set-TextValue "D31" "3.01"
$ws.Range("E31").Value = "  -5.92%  "
Set-TextValue "D32" "2.24"
$ws.Range("E32").Value = "  -1.28%  "
$ws.Range("E33").Value = "  -2.47%  "
Set-TextValue "D34" "29.18"
$ws.Range("E34").Value = "  -1.00%  "
$ws.Range("E35").Value = "  -3.83%  "
Set-TextValue "D36" "9.21"
$ws.Range("E36").Value = "  +0.20%  "
Set-TextValue "D37" "3.775.03"
$ws.Range("E37").Value = "  -0.82%  "
Set-TextValue "D38" "0.999"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  -2.13%  "
Set-TextValue "D40" "6.38"
$ws.Range("E40").Value = "  +5.61%  "
Set-TextValue "D41" "2.43"
$ws.Range("E41").Value = "  +2.79%  "
$ws.Range("E42").Value = "  -1.65%  "
$ws.Range("E43").Value = "  -2.32%  "
Set-TextValue "D45" "172.69"
$ws.Range("E45").Value = "  +5.50%  "
$ws.Range("E46").Value = "  +0.17%  "
Set-TextValue "D49" "426.22"
$ws.Range("E49").Value = "  +2.59%  "
$ws.Range("E50").Value = "  -1.32%  "
Set-TextValue "D51" "8.59"
$ws.Range("E51").Value = "  -0.34%  "

# Row 47/48: swap FLOKI and OKB
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D47" "49.90"
$ws.Range("E47").Value = "  +2.58%  "

$ws.Range("B48").Value = "FLOKI"
$ws.Range("C48").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextValue "D48" "0.000308"
$ws.Range("E48").Value = "  -3.10%  "
